$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new value looks like a plain decimal number must be
# forced to Text format first, otherwise Excel silently stores them as a
# floating point number and the original text (e.g. trailing zeros) is lost.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply updated cell values row by row, in sheet order.
$ws.Range('D2').Value = '46.977.40'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '2.478.11'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '319.40'
$ws.Range('E5').Value = '  -1.21%  '
$ws.Range('D6').Value = '107.62'
$ws.Range('E6').Value = '  +2.47%  '
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -1.33%  '
$ws.Range('D10').Value = '38.68'
$ws.Range('E10').Value = '  +7.24%  '
$ws.Range('E11').Value = '  -1.08%  '
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = '18.08'
$ws.Range('E13').Value = '  -0.87%  '
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').Value = '2.864.51'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').Value = '2.470.27'
$ws.Range('E16').Value = '  -2.75%  '
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').Value = '46.926.44'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').Value = '12.65'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('D21').Value = '2.75'
$ws.Range('E21').Value = '  +14.41%  '
$ws.Range('E22').Value = '  -0.76%  '
$ws.Range('D23').Value = '70.25'
$ws.Range('E23').Value = '  -0.46%  '
$ws.Range('D24').Value = '244.76'
$ws.Range('E24').Value = '  -1.73%  '
$ws.Range('D25').Value = '2.55'
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '25.51'
$ws.Range('E27').Value = '  -2.48%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '10.00'
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('E30').Value = '  +4.59%  '
$ws.Range('D31').Value = '34.80'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').Value = '49.41'
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('D33').Value = '19.76'
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  +1.14%  '
$ws.Range('D36').Value = '1.01'
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('E37').Value = '  +2.11%  '
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('E39').Value = '  -0.57%  '
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').Value = '118.37'
$ws.Range('E42').Value = '  -3.94%  '
$ws.Range('D43').Value = '21.57'
$ws.Range('E43').Value = '  +4.60%  '
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('D45').Value = '1.973.17'
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('E46').Value = '  +0.78%  '
$ws.Range('E47').Value = '  -5.27%  '
$ws.Range('D48').Value = '9.02'
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('E49').Value = '  -3.09%  '
$ws.Range('E50').Value = '  -4.57%  '
$ws.Range('D51').Value = '56.75'
$ws.Range('E51').Value = '  +4.31%  '
